$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 - new entry dated 13.03.2020
$ws.Range("B11").Value = 43903
$ws.Range("C11").Value = 0.5625
$ws.Range("D11").Value = 0.62847222222222221
$ws.Range("F11").Value = "CLion"
$ws.Range("G11").Value = "Bataille Navale"
$ws.Range("H11").Value = "Programmation du jeu"
$ws.Range("I11").Value = "Création des différentes fonctions"

# Row 12 - new entry dated 13.03.2020 (no end time yet)
$ws.Range("B12").Value = 43903
$ws.Range("C12").Value = 0.86458333333333337
$ws.Range("F12").Value = "CLion"
$ws.Range("G12").Value = "Bataille Navale"
$ws.Range("H12").Value = "Programmation du jeu"
$ws.Range("I12").Value = "Création des différentes fonctions"

# Update active selection to match the recorded state
$ws.Range("D12").Select()
